# Update "想去人数" (column F) figures across all four sheets to match the
# newly generated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 268
$ws.Range("F6").Value = 85
$ws.Range("F7").Value = 831
$ws.Range("F8").Value = 20
$ws.Range("F9").Value = 482
$ws.Range("F10").Value = 68
$ws.Range("F11").Value = 283
$ws.Range("F14").Value = 20
$ws.Range("F15").Value = 402
$ws.Range("F16").Value = 6502
$ws.Range("F18").Value = 65
$ws.Range("F20").Value = 7435
$ws.Range("F21").Value = 33
$ws.Range("F23").Value = 3353
$ws.Range("F24").Value = 18
$ws.Range("F25").Value = 1114
$ws.Range("F26").Value = 867
$ws.Range("F28").Value = 344
$ws.Range("F29").Value = 58
$ws.Range("F30").Value = 189
$ws.Range("F31").Value = 179
$ws.Range("F32").Value = 1501
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 141
$ws.Range("F37").Value = 1128
$ws.Range("F38").Value = 1623
$ws.Range("F39").Value = 2115

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 64
$ws.Range("F4").Value = 45

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 239
$ws.Range("F3").Value = 1202

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 239
$ws.Range("F4").Value = 1202
$ws.Range("F8").Value = 268
$ws.Range("F9").Value = 85
$ws.Range("F10").Value = 831
$ws.Range("F11").Value = 20
$ws.Range("F12").Value = 482
$ws.Range("F14").Value = 68
$ws.Range("F15").Value = 283
$ws.Range("F16").Value = 64
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 402
$ws.Range("F21").Value = 6502
$ws.Range("F23").Value = 65
$ws.Range("F25").Value = 7435
$ws.Range("F26").Value = 33
$ws.Range("F28").Value = 3353
$ws.Range("F29").Value = 18
$ws.Range("F30").Value = 1114
$ws.Range("F31").Value = 867
$ws.Range("F33").Value = 344
$ws.Range("F34").Value = 58
$ws.Range("F35").Value = 45
$ws.Range("F36").Value = 189
$ws.Range("F37").Value = 179
$ws.Range("F38").Value = 1501
$ws.Range("F39").Value = 2
$ws.Range("F40").Value = 141
$ws.Range("F43").Value = 1129
$ws.Range("F44").Value = 1623
$ws.Range("F46").Value = 2115
